$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description for the "upside risk1" test row: it was testing
# full observations, now it documents the subset-observations test instead.
$ws.Range("B106").Value = "Test upside potential for subset observations"

# Add the new downside_risk macro test rows (109-111)
$ws.Range("A109").Value = "downside risk1"
$ws.Range("B109").Value = "Test downside potential for full observations"
$ws.Range("C109").Value = "downside_risk_test1"

$ws.Range("A110").Value = "downside risk2"
$ws.Range("A111").Value = "downside risk3"

$ws.Range("B110").Value = "Test downside risk for full observations"
$ws.Range("C110").Value = "downside_risk_test2"

$ws.Range("B111").Value = "Test downside risk for subset observations"
$ws.Range("C111").Value = "downside_risk_test3"

# Reposition the view/selection to match the author's final cursor position
$win = $excel.ActiveWindow
$win.ScrollRow = 91
$win.ScrollColumn = 1
$ws.Range("C118").Select() | Out-Null
